$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.496.82"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "2.069.27"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.53"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.47"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0780"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.06%  "

$ws.Range("E11").Value = "  +1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.90"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").Value = "2.374.19"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.93"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.762"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "2.061.74"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "37.424.57"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.53"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("E20").Value = "  -2.59%  "

$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.94"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.66"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +7.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.61"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.90%  "

$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.46"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.61"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0634"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.63"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("E40").Value = "  +7.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.48"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.20"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +4.76%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.98%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0954"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.76%  "

$ws.Range("D45").Value = "1.479.70"
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.68"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.27"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.95"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -5.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("D51").Value = "2.259.11"
$ws.Range("E51").Value = "  -0.09%  "
